# Fruta / hortaliza, semanal
# The edit reshuffles the (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg) tuple among the existing data
# rows (2-34). Row 17 keeps its own values. We snapshot the original
# values first, then write them back out to their new destination rows
# according to the mapping below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (values that used to live in sourceRow now
# live in destinationRow)
$mapping = @{
    2  = 12
    3  = 22
    4  = 28
    5  = 31
    6  = 25
    7  = 10
    8  = 3
    9  = 7
    10 = 14
    11 = 4
    12 = 19
    13 = 2
    14 = 20
    15 = 18
    16 = 27
    17 = 17
    18 = 30
    19 = 8
    20 = 11
    21 = 13
    22 = 26
    23 = 15
    24 = 9
    25 = 33
    26 = 5
    27 = 24
    28 = 6
    29 = 34
    30 = 16
    31 = 29
    32 = 21
    33 = 23
    34 = 32
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the original values of the affected columns for every row
# before any writes happen, since the mapping permutes values across
# rows and later reads must not see already-overwritten data.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    if ($srcRow -eq $destRow) {
        continue
    }
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcValues[$col]
    }
}
